# Auto-generated Excel COM-interop script
# Applies the numeric cell updates from the commit's target diff across all affected worksheets.
# Values were reverse-engineered cell-by-cell from the unified OOXML diff, grouped per sheet tab.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 221.33333
$ws.Range("I58").Value = 221.33333
$ws.Range("K58").Value = 663.99999
$ws.Range("M58").Value = -513.99999
$ws.Range("H64").Value = 8000
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 9666.666999999999
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 9666.666999999999
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -10162.667
$ws.Range("H67").Value = 8000
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 9666.666999999999
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 9666.666999999999
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -11382.667
$ws.Range("H107").Value = 2690
$ws.Range("I107").Value = 2613.125
$ws.Range("K107").Value = 2613.125
$ws.Range("M107").Value = -693.125
$ws.Range("H133").Value = 65999.2
$ws.Range("J133").Value = 65999.2
$ws.Range("L133").Value = 65999.2
$ws.Range("N133").Value = -76119.2
$ws.Range("H136").Value = 69374.125
$ws.Range("J136").Value = 69374.125
$ws.Range("L136").Value = 69374.125
$ws.Range("N136").Value = -79574.125
$ws.Range("H138").Value = 3766.4167
$ws.Range("J138").Value = 3609.926
$ws.Range("L138").Value = 10829.778
$ws.Range("N138").Value = -21109.778

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3631.25
$ws.Range("I61").Value = 2931
$ws.Range("K61").Value = 2931
$ws.Range("M61").Value = -2719
$ws.Range("H74").Value = 12821890
$ws.Range("I74").Value = 17545100
$ws.Range("K74").Value = 17545100
$ws.Range("M74").Value = -17544226
$ws.Range("H77").Value = 12821890
$ws.Range("I77").Value = 17545100
$ws.Range("K77").Value = 87725500
$ws.Range("M77").Value = -87721132
$ws.Range("H94").Value = 10594.8
$ws.Range("J94").Value = 10594.8
$ws.Range("L94").Value = 10594.8
$ws.Range("N94").Value = -12396.8
$ws.Range("H136").Value = 3631.25
$ws.Range("I136").Value = 2931
$ws.Range("K136").Value = 8793
$ws.Range("M136").Value = -6243

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6277.091
$ws.Range("I134").Value = 3074.5715
$ws.Range("J134").Value = 11881.5
$ws.Range("K134").Value = 9223.7145
$ws.Range("L134").Value = 35644.5
$ws.Range("M134").Value = -6688.7145
$ws.Range("N134").Value = -40714.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1049.5
$ws.Range("I22").Value = 199.44444
$ws.Range("J22").Value = 3599.6667
$ws.Range("K22").Value = 199.44444
$ws.Range("L22").Value = 3599.6667
$ws.Range("M22").Value = 150.55556
$ws.Range("N22").Value = -4299.6667
$ws.Range("H31").Value = 35462.97
$ws.Range("I31").Value = 3618.889
$ws.Range("J31").Value = 71287.56
$ws.Range("K31").Value = 3618.889
$ws.Range("L31").Value = 71287.56
$ws.Range("M31").Value = -3323.889
$ws.Range("N31").Value = -71877.56
$ws.Range("H34").Value = 35462.97
$ws.Range("I34").Value = 3618.889
$ws.Range("J34").Value = 71287.56
$ws.Range("K34").Value = 3618.889
$ws.Range("L34").Value = 71287.56
$ws.Range("M34").Value = -3416.889
$ws.Range("N34").Value = -71691.56
$ws.Range("H62").Value = 11667.667
$ws.Range("I62").Value = 2998
$ws.Range("K62").Value = 2998
$ws.Range("M62").Value = -2374
$ws.Range("H65").Value = 11667.667
$ws.Range("I65").Value = 2998
$ws.Range("K65").Value = 14990
$ws.Range("M65").Value = -11870
$ws.Range("H132").Value = 4322.3955
$ws.Range("I132").Value = 4280.3784
$ws.Range("J132").Value = 4581.5
$ws.Range("K132").Value = 12841.1352
$ws.Range("L132").Value = 13744.5
$ws.Range("M132").Value = -10311.1352
$ws.Range("N132").Value = -18804.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2400.8
$ws.Range("J5").Value = 9103.666999999999
$ws.Range("L5").Value = 27311.001
$ws.Range("N5").Value = -27535.001
$ws.Range("H8").Value = 403.85715
$ws.Range("I8").Value = 403.85715
$ws.Range("K8").Value = 1211.57145
$ws.Range("M8").Value = -1072.57145
$ws.Range("H135").Value = 2400.8
$ws.Range("J135").Value = 9103.666999999999
$ws.Range("L135").Value = 81933.003
$ws.Range("N135").Value = -87003.003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H80").Value = 7217
$ws.Range("I80").Value = 1565
$ws.Range("J80").Value = 9101
$ws.Range("K80").Value = 1565
$ws.Range("L80").Value = 9101
$ws.Range("M80").Value = -567
$ws.Range("N80").Value = -11097
$ws.Range("H83").Value = 7217
$ws.Range("I83").Value = 1565
$ws.Range("J83").Value = 9101
$ws.Range("K83").Value = 7825
$ws.Range("L83").Value = 45505
$ws.Range("M83").Value = -2833
$ws.Range("N83").Value = -55489
$ws.Range("H123").Value = 38364.7
$ws.Range("J123").Value = 38364.7
$ws.Range("L123").Value = 38364.7
$ws.Range("N123").Value = -43264.7
$ws.Range("H135").Value = 64999.25
$ws.Range("J135").Value = 64999.25
$ws.Range("L135").Value = 64999.25
$ws.Range("N135").Value = -75139.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11340
$ws.Range("I22").Value = 3500
$ws.Range("J22").Value = 13300
$ws.Range("K22").Value = 3500
$ws.Range("L22").Value = 13300
$ws.Range("M22").Value = -3205
$ws.Range("N22").Value = -13890
$ws.Range("H27").Value = 11340
$ws.Range("I27").Value = 3500
$ws.Range("J27").Value = 13300
$ws.Range("K27").Value = 3500
$ws.Range("L27").Value = 13300
$ws.Range("M27").Value = -3393
$ws.Range("N27").Value = -13514
$ws.Range("H46").Value = 4946.154
$ws.Range("I46").Value = 3499.75
$ws.Range("K46").Value = 3499.75
$ws.Range("M46").Value = -3311.75
$ws.Range("H55").Value = 1257.4054
$ws.Range("I55").Value = 712.7778
$ws.Range("J55").Value = 1773.3684
$ws.Range("K55").Value = 712.7778
$ws.Range("L55").Value = 1773.3684
$ws.Range("M55").Value = -539.7778
$ws.Range("N55").Value = -2119.3684
$ws.Range("H64").Value = 26333.334
$ws.Range("J64").Value = 34500
$ws.Range("L64").Value = 34500
$ws.Range("N64").Value = -34950
$ws.Range("H67").Value = 26333.334
$ws.Range("J67").Value = 34500
$ws.Range("L67").Value = 34500
$ws.Range("N67").Value = -36060
$ws.Range("H74").Value = 60000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 60000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 60000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -61996
$ws.Range("H77").Value = 60000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 60000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 180000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -189984
$ws.Range("H136").Value = 4159.0454
$ws.Range("I136").Value = 2182.9333
$ws.Range("J136").Value = 8393.571
$ws.Range("K136").Value = 6548.7999
$ws.Range("L136").Value = 25180.713
$ws.Range("M136").Value = -3998.7999
$ws.Range("N136").Value = -30280.713
$ws.Range("H139").Value = 69749
$ws.Range("J139").Value = 69749
$ws.Range("L139").Value = 69749
$ws.Range("N139").Value = -80029

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 49984.332
$ws.Range("J93").Value = 49984.332
$ws.Range("L93").Value = 49984.332
$ws.Range("N93").Value = -54976.332
$ws.Range("H107").Value = 620.6875
$ws.Range("I107").Value = 638.2857
$ws.Range("K107").Value = 1914.8571
$ws.Range("M107").Value = 5.142899999999827
$ws.Range("H132").Value = 4408.6
$ws.Range("I132").Value = 4250.56
$ws.Range("J132").Value = 5198.8
$ws.Range("K132").Value = 12751.68
$ws.Range("L132").Value = 15596.4
$ws.Range("M132").Value = -10221.68
$ws.Range("N132").Value = -20656.4
$ws.Range("H136").Value = 3076.4062
$ws.Range("I136").Value = 1253.3704
$ws.Range("K136").Value = 3760.1112
$ws.Range("M136").Value = -1210.1112
$ws.Range("H140").Value = 60624.75
$ws.Range("J140").Value = 60624.75
$ws.Range("L140").Value = 60624.75
$ws.Range("N140").Value = -70984.75
$ws.Range("H141").Value = 96535.5
$ws.Range("J141").Value = 96535.5
$ws.Range("L141").Value = 96535.5
$ws.Range("N141").Value = -106895.5
